$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '27.859.66'
Set-TextValue 'D3' '1.872.07'
Set-TextValue 'D5' '313.08'
Set-TextValue 'E5' '  +0.58%  '
Set-TextValue 'E6' '  -0.64%  '
Set-TextValue 'E7' '  +0.70%  '
Set-TextValue 'D8' '0.3820'
Set-TextValue 'E8' '  +2.90%  '
Set-TextValue 'D9' '0.07372'
Set-TextValue 'E9' '  +1.16%  '
Set-TextValue 'D10' '0.9399'
Set-TextValue 'E10' '  +0.59%  '
Set-TextValue 'D11' '21.03'
Set-TextValue 'E11' '  +4.69%  '
Set-TextValue 'D12' '0.07795'
Set-TextValue 'E12' '  -0.83%  '
Set-TextValue 'D13' '1.915.04'
Set-TextValue 'E13' '  +2.43%  '
Set-TextValue 'D14' '5.508'
Set-TextValue 'D15' '6.614'
Set-TextValue 'E15' '  +1.19%  '
Set-TextValue 'D16' '91.36'
Set-TextValue 'E16' '  +1.55%  '
Set-TextValue 'D17' '1.013'
Set-TextValue 'E17' '  -0.73%  '
Set-TextValue 'D18' '0.000008865'
Set-TextValue 'E18' '  +1.66%  '
Set-TextValue 'E19' '  -0.66%  '
Set-TextValue 'D20' '27.883.91'
Set-TextValue 'E20' '  +2.40%  '
Set-TextValue 'E21' '  +1.25%  '
Set-TextValue 'E22' '  +0.54%  '
Set-TextValue 'D23' '2.136.68'
Set-TextValue 'E23' '  +1.64%  '
Set-TextValue 'E24' '  +1.78%  '
Set-TextValue 'D25' '157.68'
Set-TextValue 'E25' '  +2.68%  '
Set-TextValue 'D26' '1.946'
Set-TextValue 'E26' '  -0.24%  '
Set-TextValue 'E27' '  +0.61%  '
Set-TextValue 'D28' '2.048'
Set-TextValue 'E28' '  +2.60%  '
Set-TextValue 'D29' '116.02'
Set-TextValue 'D30' '4.980'
Set-TextValue 'E30' '  +0.99%  '
Set-TextValue 'D31' '0.08899'
Set-TextValue 'E31' '  +0.15%  '
Set-TextValue 'E32' '  +0.67%  '
Set-TextValue 'D33' '1.226'
Set-TextValue 'E33' '  +3.62%  '
Set-TextValue 'D34' '0.7718'
Set-TextValue 'E34' '  +4.76%  '
Set-TextValue 'D35' '4.657'
Set-TextValue 'E35' '  +1.76%  '
Set-TextValue 'D36' '2.734'
Set-TextValue 'E36' '  +1.62%  '
Set-TextValue 'D37' '1.132'
Set-TextValue 'E37' '  +0.90%  '
Set-TextValue 'E38' '  +1.85%  '
Set-TextValue 'D39' '0.5621'
Set-TextValue 'E39' '  +5.54%  '
Set-TextValue 'D40' '0.05385'
Set-TextValue 'E40' '  +2.52%  '
Set-TextValue 'E41' '  +0.08%  '
Set-TextValue 'D42' '7.047'
Set-TextValue 'E42' '  -0.24%  '
Set-TextValue 'D43' '8.551'
Set-TextValue 'D44' '0.1527'
Set-TextValue 'E44' '  -0.06%  '
Set-TextValue 'E45' '  +1.34%  '
Set-TextValue 'D46' '0.4876'
Set-TextValue 'E46' '  +2.22%  '
Set-TextValue 'D47' '105.46'
Set-TextValue 'E47' '  +2.94%  '
Set-TextValue 'E48' '  -0.67%  '
Set-TextValue 'D49' '1.663'
Set-TextValue 'E49' '  +2.12%  '
Set-TextValue 'D50' '68.12'
Set-TextValue 'E50' '  +2.58%  '
Set-TextValue 'D51' '0.06121'
Set-TextValue 'E51' '  +0.75%  '
